$d = $word.ActiveDocument

$replacements = @(
    @("976÷7=139, 3", "510÷7=72, 6"),
    @("408÷3=136, 0", "134÷8=16, 6"),
    @("217÷4=54, 1", "930÷5=186, 0"),
    @("722÷5=144, 2", "205÷3=68, 1"),
    @("628÷6=104, 4", "109÷4=27, 1"),
    @("968÷7=138, 2", "753÷6=125, 3"),
    @("674÷4=168, 2", "848÷3=282, 2"),
    @("194÷9=21, 5", "340÷7=48, 4"),
    @("651÷9=72, 3", "369÷6=61, 3"),
    @("177÷5=35, 2", "950÷8=118, 6"),
    @("724÷5=144, 4", "793÷7=113, 2"),
    @("864÷9=96, 0", "538÷3=179, 1"),
    @("838÷7=119, 5", "928÷6=154, 4"),
    @("539÷3=179, 2", "663÷6=110, 3"),
    @("854÷9=94, 8", "415÷2=207, 1"),
    @("122÷3=40, 2", "712÷4=178, 0"),
    @("603÷9=67, 0", "744÷8=93, 0"),
    @("464÷3=154, 2", "536÷9=59, 5"),
    @("380÷2=190, 0", "595÷2=297, 1"),
    @("818÷9=90, 8", "141÷6=23, 3"),
    @("100÷7=14, 2", "711÷7=101, 4"),
    @("888÷5=177, 3", "879÷7=125, 4"),
    @("176÷9=19, 5", "855÷2=427, 1"),
    @("234÷2=117, 0", "245÷7=35, 0"),
    @("109÷7=15, 4", "742÷4=185, 2")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
